$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) values for B:E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 values for B:E
$ws.Range("B2").Value = 517.64724732000002
$ws.Range("C2").Value = 435.62734945312508
$ws.Range("D2").Value = 518.09924965750008
$ws.Range("E2").Value = 434.55029483812501

# Update row 3 values for B:E
$ws.Range("B3").Value = 522.69955899999991
$ws.Range("C3").Value = 431.21804616000003
$ws.Range("D3").Value = 527.48455583999998
$ws.Range("E3").Value = 439.87016094749998

# Update the active selection to B1:E3
$ws.Range("B1:E3").Select()
